$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): existing A1:H1 get new text, I1:N1 are brand new columns ---
$ws.Range("A1").Value = "Segment Name"
$ws.Range("B1").Value = "Segment Description"
$ws.Range("C1").Value = "Segment Status"
$ws.Range("D1").Value = "Segment Lifetime"
$ws.Range("E1").Value = "Trait Folder Path"
$ws.Range("F1").Value = "Data Source ID"
$ws.Range("G1").Value = "Data Source Name"
$ws.Range("H1").Value = "Data Source Result"
$ws.Range("I1").Value = "Data Feed Result"
$ws.Range("J1").Value = "Segments and Overlap Plan Result"
$ws.Range("K1").Value = "Modeling Plan Result"
$ws.Range("L1").Value = "Activation Plan Result"
$ws.Range("M1").Value = "Trait Folder Result"
$ws.Range("N1").Value = "Create Segment Result"

# Give the newly added header cells (I1:N1) the same bold/bordered style as the
# rest of the header row so the whole row reuses a single cell style.
$ws.Range("A1").Copy()
$ws.Range("I1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$description = "HG Data's curation IP and supervised machine learning utilizes billions of unstructured source documents for unseen intelligence and insights about industries and companies worldwide. Use this segment to address the key decision makers and influencers with detected past purchases of the category of technologies listed here."

# --- Data rows (2 and 3): populate column-by-column ---
$ws.Range("A2").Value = "Test Segment 1"
$ws.Range("A3").Value = "Test Segment 2"

$ws.Range("B2").Value = $description
$ws.Range("B3").Value = $description

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("D2").Value = 90
$ws.Range("D3").Value = 90

$ws.Range("E2").Value = "/All Traits/TEST20181030/TEST"
$ws.Range("E3").Value = "/All Traits/TEST20181030"

$ws.Range("F2").Value = 421167
$ws.Range("F3").Value = 421167

$ws.Range("G2").Value = "test20181030"
$ws.Range("G3").Value = "test20181030"

$ws.Range("H2").Value = "Existing Data Source"
$ws.Range("H3").Value = "Existing Data Source"

$ws.Range("I2").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("L3").ClearContents()

$ws.Range("M2").Value = "Created"
$ws.Range("M3").Value = "Existing folder"

$ws.Range("N2").Value = "['description size must be between 0 and 255']"
$ws.Range("N3").Value = "['description size must be between 0 and 255']"
